$d = $word.ActiveDocument

# The "Bibliografia" paragraph originally has every reference concatenated
# into a single run of text. We split it into multiple <w:t> runs separated
# by manual line breaks (<w:br/>) by inserting "^l" (Word's Find/Replace
# code for a manual line break) at each reference boundary.

$pairs = @(
    @("Cengage, 2ed, 2022.Nilo Ney", "Cengage, 2ed, 2022.^lNilo Ney"),
    @("3a ed, 2019.Ramalho, L.", "3a ed, 2019.^lRamalho, L."),
    @("O’Reilly-Novatec, 2015Downey, A. B.", "O’Reilly-Novatec, 2015^lDowney, A. B."),
    @("O’Reilly-Novatec, 2016.STEWART,", "O’Reilly-Novatec, 2016.^lSTEWART,"),
    @("Cambridge University Press, 2014.TELLES,", "Cambridge University Press, 2014.^lTELLES,"),
    @("Thomson Course Technology PTR, 2008.LUTZ,", "Thomson Course Technology PTR, 2008.^lLUTZ,"),
    @("Sebastopol, CA: O’Reilly Media, 2006.MCGREGGOR,", "Sebastopol, CA: O’Reilly Media, 2006.^lMCGREGGOR,")
)

foreach ($pair in $pairs) {
    $find = $pair[0]
    $replace = $pair[1]
    $range = $d.Content
    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}
